$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column B, filled in order 18, 17, 19 (matches authoring order of shared strings)
$ws.Range("B18").Value = "2) Assigning team member roles to design work flow diagram"
$ws.Range("B17").Value = "1) Pitching project idea, discussed pros and cons of project and assigned team members"
$ws.Range("B19").Value = "3)Assigning team member roles is bit difficul"

# Column C, filled top to bottom
$ws.Range("C17").Value = "1)We as a team discussed about software methodology to follow and discussed about the number of layouts to be placed and software requirements to be installed on every team members device."
$ws.Range("C18").Value = "2)Getting knowledge and designing basic layouts of the application"
$ws.Range("C19").Value = "3)Integrating the frontend and database"

# Column D, filled top to bottom
$ws.Range("D17").Value = "1)Adithya and me shared roles as frontend developers and we are anticipated 5 layout in the application and we discussed with team members about layouts and how layouts has to be designed"
$ws.Range("D18").Value = "2)Concentrating on UI and layouts"
$ws.Range("D19").Value = "3)Finalizing UI design"

# Row height updates
$ws.Rows.Item(17).RowHeight = 93
$ws.Rows.Item(18).RowHeight = 31
$ws.Rows.Item(19).RowHeight = 31.5

# Update sheet view - scroll position and active cell selection
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("D19").Select()
